# Scheduled runner update: refresh market-board derived values across the
# per-job profit sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# ALC
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("ALC")

$ws.Range("H17").Value = 308491.38
$ws.Range("J17").Value = 308491.38
$ws.Range("L17").Value = 925474.14
$ws.Range("N17").Value = -925810.14

$ws.Range("H19").Value = 1214.9778
$ws.Range("J19").Value = 1353.75
$ws.Range("L19").Value = 1353.75
$ws.Range("N19").Value = -1703.75

$ws.Range("H64").Value = 5725.7144
$ws.Range("J64").Value = 3928.889
$ws.Range("L64").Value = 3928.889
$ws.Range("N64").Value = -4424.889

$ws.Range("H67").Value = 5725.7144
$ws.Range("J67").Value = 3928.889
$ws.Range("L67").Value = 3928.889
$ws.Range("N67").Value = -5644.889

$ws.Range("H125").Value = 19136.666
$ws.Range("I125").Value = 9805
$ws.Range("J125").Value = 37800
$ws.Range("K125").Value = 88245
$ws.Range("L125").Value = 340200
$ws.Range("M125").Value = -85785
$ws.Range("N125").Value = -345120

$ws.Range("H137").Value = 1146.225
$ws.Range("I137").Value = 916.14703
$ws.Range("K137").Value = 2748.44109
$ws.Range("M137").Value = -198.4410899999998

$ws.Range("H138").Value = 3490.09
$ws.Range("I138").Value = 679.03125
$ws.Range("J138").Value = 4812.9414
$ws.Range("K138").Value = 2037.09375
$ws.Range("L138").Value = 14438.8242
$ws.Range("M138").Value = 3102.90625
$ws.Range("N138").Value = -24718.8242

# ---------------------------------------------------------------------
# ARM
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("ARM")

$ws.Range("H32").Value = 4584697
$ws.Range("I32").Value = 4920939.5
$ws.Range("K32").Value = 4920939.5
$ws.Range("M32").Value = -4920652.5

$ws.Range("H88").Value = 4412.7
$ws.Range("I88").Value = 3753
$ws.Range("J88").Value = 4577.625
$ws.Range("K88").Value = 3753
$ws.Range("L88").Value = 4577.625
$ws.Range("M88").Value = -3347
$ws.Range("N88").Value = -5389.625

$ws.Range("H91").Value = 4412.7
$ws.Range("I91").Value = 3753
$ws.Range("J91").Value = 4577.625
$ws.Range("K91").Value = 3753
$ws.Range("L91").Value = 4577.625
$ws.Range("M91").Value = -2349
$ws.Range("N91").Value = -7385.625

# ---------------------------------------------------------------------
# BSM
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("BSM")

$ws.Range("H13").Value = 55000
$ws.Range("J13").Value = 55000
$ws.Range("L13").Value = 55000
$ws.Range("N13").Value = -55336

$ws.Range("H86").Value = 1953.4688
$ws.Range("I86").Value = 1827.0769
$ws.Range("J86").Value = 2501.1667
$ws.Range("K86").Value = 1827.0769
$ws.Range("L86").Value = 2501.1667
$ws.Range("M86").Value = -704.0769
$ws.Range("N86").Value = -4747.1667

$ws.Range("H89").Value = 1953.4688
$ws.Range("I89").Value = 1827.0769
$ws.Range("J89").Value = 2501.1667
$ws.Range("K89").Value = 9135.3845
$ws.Range("L89").Value = 12505.8335
$ws.Range("M89").Value = -3519.3845
$ws.Range("N89").Value = -23737.8335

# ---------------------------------------------------------------------
# CRP
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("CRP")

$ws.Range("H31").Value = 2910.9148
$ws.Range("I31").Value = 3002.9167
$ws.Range("J31").Value = 2609.818
$ws.Range("K31").Value = 3002.9167
$ws.Range("L31").Value = 2609.818
$ws.Range("M31").Value = -2707.9167
$ws.Range("N31").Value = -3199.818

$ws.Range("H34").Value = 2910.9148
$ws.Range("I34").Value = 3002.9167
$ws.Range("J34").Value = 2609.818
$ws.Range("K34").Value = 3002.9167
$ws.Range("L34").Value = 2609.818
$ws.Range("M34").Value = -2800.9167
$ws.Range("N34").Value = -3013.818

$ws.Range("H132").Value = 1427.5
$ws.Range("I132").Value = 918.88
$ws.Range("K132").Value = 2756.64
$ws.Range("M132").Value = -226.6399999999999

$ws.Range("H134").Value = 1448.9131
$ws.Range("I134").Value = 592.8333
$ws.Range("J134").Value = 2382.818
$ws.Range("K134").Value = 1778.4999
$ws.Range("L134").Value = 7148.454000000001
$ws.Range("M134").Value = 756.5001
$ws.Range("N134").Value = -12218.454

# ---------------------------------------------------------------------
# CUL
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("CUL")

$ws.Range("H105").Value = 156142.86
$ws.Range("J105").Value = 156142.86
$ws.Range("L105").Value = 468428.58
$ws.Range("N105").Value = -473670.58

$ws.Range("H126").Value = 47858.957
$ws.Range("I126").Value = 168898.33
$ws.Range("J126").Value = 5139.1763
$ws.Range("K126").Value = 506694.99
$ws.Range("L126").Value = 15417.5289
$ws.Range("M126").Value = -501754.99
$ws.Range("N126").Value = -25297.5289

# ---------------------------------------------------------------------
# GSM
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("GSM")

$ws.Range("H70").Value = 4309.3125
$ws.Range("I70").Value = 4258.9
$ws.Range("J70").Value = 4393.3335
$ws.Range("K70").Value = 4258.9
$ws.Range("L70").Value = 4393.3335
$ws.Range("M70").Value = -3988.9
$ws.Range("N70").Value = -4933.3335

$ws.Range("H73").Value = 4309.3125
$ws.Range("I73").Value = 4258.9
$ws.Range("J73").Value = 4393.3335
$ws.Range("K73").Value = 4258.9
$ws.Range("L73").Value = 4393.3335
$ws.Range("M73").Value = -3322.9
$ws.Range("N73").Value = -6265.3335

# ---------------------------------------------------------------------
# LTW
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("LTW")

$ws.Range("H82").Value = 4549.75
$ws.Range("I82").Value = 3899.1667
$ws.Range("J82").Value = 6501.5
$ws.Range("K82").Value = 3899.1667
$ws.Range("L82").Value = 6501.5
$ws.Range("M82").Value = -3538.1667
$ws.Range("N82").Value = -7223.5

$ws.Range("H85").Value = 4549.75
$ws.Range("I85").Value = 3899.1667
$ws.Range("J85").Value = 6501.5
$ws.Range("K85").Value = 3899.1667
$ws.Range("L85").Value = 6501.5
$ws.Range("M85").Value = -2651.1667
$ws.Range("N85").Value = -8997.5

$ws.Range("H122").Value = 100000
$ws.Range("J122").Value = 0
$ws.Range("L122").Value = 0
$ws.Range("N122").ClearContents()

$ws.Range("H132").Value = 1696.8422
$ws.Range("I132").Value = 1494.6538
$ws.Range("J132").Value = 3799.6
$ws.Range("K132").Value = 4483.9614
$ws.Range("L132").Value = 11398.8
$ws.Range("M132").Value = -1953.9614
$ws.Range("N132").Value = -16458.8

# ---------------------------------------------------------------------
# WVR
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("WVR")

$ws.Range("H113").Value = 246.5
$ws.Range("I113").Value = 247.26666
$ws.Range("J113").Value = 244.2
$ws.Range("K113").Value = 741.79998
$ws.Range("L113").Value = 732.5999999999999
$ws.Range("M113").Value = 1428.20002
$ws.Range("N113").Value = -5072.6

$ws.Range("H122").Value = 551.5833
$ws.Range("I122").Value = 440.44446
$ws.Range("J122").Value = 885
$ws.Range("K122").Value = 1321.33338
$ws.Range("L122").Value = 2655
$ws.Range("M122").Value = 1128.66662
$ws.Range("N122").Value = -7555

$ws.Range("H123").Value = 21401.875
$ws.Range("J123").Value = 21401.875
$ws.Range("L123").Value = 21401.875
$ws.Range("N123").Value = -31201.875
